$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 333, shifting existing rows 333:401 down to 334:402
$ws.Rows.Item(333).Insert()

# Populate the newly inserted row 333 with the new record
$ws.Cells.Item(333, 1).Value = 3
$ws.Cells.Item(333, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(333, 3).Value = "Coquimbo"
$ws.Cells.Item(333, 4).Value = 44511
$ws.Cells.Item(333, 5).Value = 5
$ws.Cells.Item(333, 6).Value = 100112024
$ws.Cells.Item(333, 7).Value = "Choclo"
$ws.Cells.Item(333, 8).Value = "Dulce o Americano"
$ws.Cells.Item(333, 9).Value = "Primera"
$ws.Cells.Item(333, 10).Value = 78
$ws.Cells.Item(333, 11).Value = 41000
$ws.Cells.Item(333, 12).Value = 42000
$ws.Cells.Item(333, 13).Value = 41513
$ws.Cells.Item(333, 14).Value = "$/malla 70 unidades"
$ws.Cells.Item(333, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(333, 16).Value = 593
$ws.Cells.Item(333, 17).Value = 70
$ws.Cells.Item(333, 18).Value = "Hortaliza"

Write-Host "done"
